$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (15) with the "Percentage revenue by year" question.
# Set the URL cell first so the shared-string table picks up the same
# insertion order as the authored workbook (URL, then title, then status).
$ws.Range("C15").Value = "https://www.interviewquery.com/questions/percentage-of-revenue-by-year"
$ws.Range("A15").Value = "Percentage revenue by year"
$ws.Range("B15").Value = "Medium"
$ws.Range("D15").Value = "Unsolved"
$ws.Range("E15").Value = ""

# Match the row styling used for the rest of the table: a themed
# (accent) font colour plus the taller 34pt row height the other
# multi-line rows use. Columns C:D already wrap via the column style.
$ws.Range("A15:E15").Font.ThemeColor = 10
$ws.Rows.Item(15).RowHeight = 34

# Move the active selection the same way the author's session ended up.
$ws.Range("C23").Select() | Out-Null
